# Update gh-pages output (广州-漫展信息.xlsx) to the newly scraped numbers.
# Sheets: 展览 (exhibitions), 演出 (shows), 本地生活 (local-life, untouched),
# 全部类型 (all-types, a date-sorted merge of the first two).

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
}

# ---------------------------------------------------------------------------
# Sheet "展览" (sheet 1): refresh "想去人数" (F) / "最低票价" (G) counters.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1Changes = @(
    @(2, 11356, $null),
    @(3, 1849, "55"),
    @(4, 513, $null),
    @(5, 812, $null),
    @(6, 2408, "已售罄"),
    @(7, 733, $null),
    @(8, 970, $null),
    @(9, 573, $null),
    @(10, 436, $null),
    @(11, 471, $null),
    @(12, 433, $null),
    @(13, 1285, $null),
    @(14, 617, $null),
    @(15, 72, $null),
    @(16, 955, $null),
    @(17, 468, $null),
    @(18, 631, $null),
    @(19, 1017, $null),
    @(20, 194, $null),
    @(21, 907, $null),
    @(22, 116, $null),
    @(23, 214, $null),
    @(24, 112, $null),
    @(25, 248, $null),
    @(26, 653, $null),
    @(27, 151, $null),
    @(28, 94, $null),
    @(29, 308, $null)
)

foreach ($chg in $ws1Changes) {
    $ws1.Cells.Item($chg[0], 6).Value2 = $chg[1]
    if ($chg[2] -ne $null) {
        Set-TextCell $ws1 $chg[0] 7 $chg[2]
    }
}

# Row 12 (U.M.A闪耀ONLY) also got a refreshed cover image.
Set-TextCell $ws1 12 10 "//i2.hdslb.com/bfs/openplatform/202401/n7B2YOEs1705908251849.png"

# ---------------------------------------------------------------------------
# Sheet "演出" (sheet 2): refresh three counters, then insert the new
# "春卷饭十周年" show as row 8 (pushing 夏川里美 down to row 9).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Cells.Item(2, 6).Value2 = 88
$ws2.Cells.Item(5, 6).Value2 = 848
$ws2.Cells.Item(7, 6).Value2 = 38

$ws2.Rows.Item(8).Insert()
# Carry the bordered/bold "index" column style down onto the new row.
$ws2.Cells.Item(7, 1).Copy()
$ws2.Cells.Item(8, 1).PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false

$ws2.Cells.Item(8, 1).Value2 = 7
Set-TextCell $ws2 8 2 "2024-03-23"
$ws2.Cells.Item(8, 3).Value2 = "广州·春卷饭 十周年  2024  专场演出"
$ws2.Cells.Item(8, 4).Value2 = "革新路124号太古仓码头54汇5号仓 太空间Livehouse"
$ws2.Cells.Item(8, 5).Value2 = "2024.03.23 20:00-03.23 22:00"
$ws2.Cells.Item(8, 6).Value2 = 248
Set-TextCell $ws2 8 7 "480"
$ws2.Cells.Item(8, 8).Value2 = $false
$ws2.Cells.Item(8, 9).Value2 = "https://show.bilibili.com/platform/detail.html?id=81186"
$ws2.Cells.Item(8, 10).Value2 = "//i1.hdslb.com/bfs/openplatform/202401/ho9rIMg21705894649801.jpeg"

# The row that shifted down (old row 8, 夏川里美) keeps its data, only the
# running index in column A needs to advance to 8.
$ws2.Cells.Item(9, 1).Value2 = 8

# ---------------------------------------------------------------------------
# Sheet "全部类型" (sheet 4): same counters as above (merged view), plus the
# new 春卷饭 show inserted in date order at row 35.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4Changes = @(
    @(2, 11356, $null),
    @(3, 1849, "55"),
    @(4, 88, $null),
    @(5, 513, $null),
    @(6, 812, $null),
    @(7, 2408, "已售罄"),
    @(8, 733, $null),
    @(9, 970, $null),
    @(11, 573, $null),
    @(12, 436, $null),
    @(13, 471, $null),
    @(14, 433, $null),
    @(15, 1285, $null),
    @(17, 617, $null),
    @(18, 72, $null),
    @(19, 848, $null),
    @(20, 955, $null),
    @(21, 468, $null),
    @(22, 631, $null),
    @(23, 1017, $null),
    @(24, 194, $null),
    @(25, 907, $null),
    @(26, 116, $null),
    @(27, 214, $null),
    @(29, 112, $null),
    @(30, 248, $null),
    @(31, 653, $null),
    @(32, 151, $null),
    @(33, 38, $null),
    @(34, 94, $null)
)

foreach ($chg in $ws4Changes) {
    $ws4.Cells.Item($chg[0], 6).Value2 = $chg[1]
    if ($chg[2] -ne $null) {
        Set-TextCell $ws4 $chg[0] 7 $chg[2]
    }
}

# Row 14 (U.M.A闪耀ONLY) cover image refresh, matching sheet 1.
Set-TextCell $ws4 14 10 "//i2.hdslb.com/bfs/openplatform/202401/n7B2YOEs1705908251849.png"

# Insert 春卷饭十周年 in its date-sorted slot (2024-03-23, right before the
# 2024-04-20 Arknights show that used to sit at row 35).
$ws4.Rows.Item(35).Insert()
$ws4.Cells.Item(34, 1).Copy()
$ws4.Cells.Item(35, 1).PasteSpecial(-4122)
$ws4.Application.CutCopyMode = $false

$ws4.Cells.Item(35, 1).Value2 = 34
Set-TextCell $ws4 35 2 "2024-03-23"
$ws4.Cells.Item(35, 3).Value2 = "广州·春卷饭 十周年  2024  专场演出"
$ws4.Cells.Item(35, 4).Value2 = "革新路124号太古仓码头54汇5号仓 太空间Livehouse"
$ws4.Cells.Item(35, 5).Value2 = "2024.03.23 20:00-03.23 22:00"
$ws4.Cells.Item(35, 6).Value2 = 248
Set-TextCell $ws4 35 7 "480"
$ws4.Cells.Item(35, 8).Value2 = $false
$ws4.Cells.Item(35, 9).Value2 = "https://show.bilibili.com/platform/detail.html?id=81186"
$ws4.Cells.Item(35, 10).Value2 = "//i1.hdslb.com/bfs/openplatform/202401/ho9rIMg21705894649801.jpeg"

# The Arknights row that shifted down to row 36 also needs its "想去人数"
# counter refreshed, and its running index bumped to 35.
$ws4.Cells.Item(36, 1).Value2 = 35
$ws4.Cells.Item(36, 6).Value2 = 308

# ...and the 夏川里美 row that shifted from 36 to 37 needs its index bumped
# to 36 (no other values changed for it).
$ws4.Cells.Item(37, 1).Value2 = 36

Write-Output "edits applied"
